# Update cryptos.xlsx price/volume data (and two row swaps) to match the
# latest GitHub Actions scrape, per commit "Updated cryptos list on Wed May  8
# 12:47:47 UTC 2024 with GitHub Actions".
#
# Column D ("Price") holds free-form text (e.g. "62.173.31", "0.0000225")
# rather than numbers, so we force the cell's number format to Text ("@")
# before writing the value -- this prevents Excel's automatic type
# detection from reinterpreting the string as a number/date and mangling
# it (e.g. turning "145.00" into 145, or "0.0000225" into scientific
# notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.173.31"
$ws.Range("E2").Value = "  -2.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.05"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.74"
$ws.Range("E5").Value = "  -1.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.00"
$ws.Range("E6").Value = "  -6.33%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -3.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.991.48"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  -6.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  -2.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  -2.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -4.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.20"
$ws.Range("E14").Value = "  -6.39%  "

$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.481.53"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.154.72"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.92"
$ws.Range("E18").Value = "  -3.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.989.89"
$ws.Range("E19").Value = "  -2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.27"
$ws.Range("E20").Value = "  -6.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.74"
$ws.Range("E21").Value = "  -4.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.675"
$ws.Range("E22").Value = "  -4.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.36"
$ws.Range("E23").Value = "  -2.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.56"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -7.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("E26").Value = "  -5.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  -5.33%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.60"
$ws.Range("E30").Value = "  -3.12%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  -6.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  -5.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.73"
$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -4.00%  "

$ws.Range("E35").Value = "  -5.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0781"
$ws.Range("E36").Value = "  -5.74%  "

$ws.Range("E37").Value = "  -6.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  -6.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.03"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.89"

$ws.Range("E41").Value = "  -12.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.112"
$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "380.59"
$ws.Range("E43").Value = "  -13.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.270"
$ws.Range("E44").Value = "  -7.05%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.748.28"
$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.03"
$ws.Range("E47").Value = "  -5.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.27"
$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.59"
$ws.Range("E51").Value = "  -6.20%  "

